$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.387.17'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '2.435.62'
$ws.Range("E3").Value = '  +0.82%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '564.61'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.34'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.71%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.44%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.110'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.32%  '
$ws.Range("B10").Value = 'TRON'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.154'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.25'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.30%  '
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.350'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.26%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000178'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.820.00'
$ws.Range("E15").Value = '  -1.12%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '62.197.29'
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '2.412.95'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '324.27'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.16'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("B23").Value = 'Litecoin'
$ws.Range("C23").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.51%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.03%  '
$ws.Range("B25").Value = 'Aptos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.80'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.33%  '
$ws.Range("B26").Value = 'Bittensor'
$ws.Range("C26").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '561.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.81%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").Value = '0.0₃0963'
$ws.Range("E27").Value = '  +1.73%  '
$ws.Range("D28").Value = '2.554.41'
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.28%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.24%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.44'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.148'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.71%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.88'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.63%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.53'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.28%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.85'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.78%  '
$ws.Range("E36").Value = '  -0.15%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.383'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.58'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("B39").Value = 'EthereumClassic'
$ws.Range("C39").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '149.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.98%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.30%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.58%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '148.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0535'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.599'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0927'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.98%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0231'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.25%  '
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '11.60'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.63%  '
